$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 517.2
$ws.Range("I2").Value = 571.5
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 571.5
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -458.5
$ws.Range("N2").Value = -526

$ws.Range("H12").Value = 897.44446
$ws.Range("I12").Value = 531.6667
$ws.Range("J12").Value = 1080.3334
$ws.Range("K12").Value = 531.6667
$ws.Range("L12").Value = 1080.3334
$ws.Range("M12").Value = -361.6667
$ws.Range("N12").Value = -1420.3334

$ws.Range("H43").Value = 1156.4
$ws.Range("I43").Value = 377.75
$ws.Range("J43").Value = 1675.5
$ws.Range("K43").Value = 377.75
$ws.Range("L43").Value = 1675.5
$ws.Range("M43").Value = -308.75
$ws.Range("N43").Value = -1813.5

$ws.Range("H45").Value = 18702768
$ws.Range("I45").Value = 51349.5
$ws.Range("J45").Value = 28028478
$ws.Range("K45").Value = 154048.5
$ws.Range("L45").Value = 84085434
$ws.Range("M45").Value = -153856.5
$ws.Range("N45").Value = -84085818

$ws.Range("H48").Value = 2500
$ws.Range("I48").Value = 2500
$ws.Range("K48").Value = 7500
$ws.Range("M48").Value = -7208

$ws.Range("H49").Value = 38461784
$ws.Range("J49").Value = 76923070
$ws.Range("L49").Value = 230769210
$ws.Range("N49").Value = -230769482

$ws.Range("H51").Value = 2053.9092
$ws.Range("I51").Value = 3373
$ws.Range("J51").Value = 1300.1428
$ws.Range("K51").Value = 3373
$ws.Range("L51").Value = 1300.1428
$ws.Range("M51").Value = -2889
$ws.Range("N51").Value = -2268.1428

$ws.Range("H53").Value = 216.44444
$ws.Range("I53").Value = 171.5
$ws.Range("J53").Value = 252.4
$ws.Range("K53").Value = 171.5
$ws.Range("L53").Value = 252.4
$ws.Range("M53").Value = 465.5
$ws.Range("N53").Value = -1526.4

$ws.Range("H56").Value = 2500
$ws.Range("I56").Value = 2500
$ws.Range("K56").Value = 7500
$ws.Range("M56").Value = -6966

$ws.Range("H59").Value = 3473107.8
$ws.Range("I59").Value = 1008.5
$ws.Range("J59").Value = 5209157.5
$ws.Range("K59").Value = 3025.5
$ws.Range("L59").Value = 15627472.5
$ws.Range("M59").Value = -2468.5
$ws.Range("N59").Value = -15628586.5

$ws.Range("H69").Value = 6538.778
$ws.Range("I69").Value = 6013
$ws.Range("J69").Value = 6959.4
$ws.Range("K69").Value = 18039
$ws.Range("L69").Value = 20878.2
$ws.Range("M69").Value = -17165
$ws.Range("N69").Value = -22626.2

$ws.Range("H72").Value = 6538.778
$ws.Range("I72").Value = 6013
$ws.Range("J72").Value = 6959.4
$ws.Range("K72").Value = 54117
$ws.Range("L72").Value = 62634.6
$ws.Range("M72").Value = -49749
$ws.Range("N72").Value = -71370.60000000001

$ws.Range("H112").Value = 1662.6207
$ws.Range("I112").Value = 490
$ws.Range("J112").Value = 1704.5
$ws.Range("K112").Value = 1470
$ws.Range("L112").Value = 5113.5
$ws.Range("M112").Value = -362
$ws.Range("N112").Value = -7329.5

$ws.Range("H129").Value = 889.57574
$ws.Range("I129").Value = 306.7143
$ws.Range("J129").Value = 958.7288
$ws.Range("K129").Value = 920.1428999999999
$ws.Range("L129").Value = 2876.1864
$ws.Range("M129").Value = 4079.8571
$ws.Range("N129").Value = -12876.1864

$ws.Range("H138").Value = 3790632.8
$ws.Range("I138").Value = 1305.9231
$ws.Range("J138").Value = 6806627.5
$ws.Range("K138").Value = 3917.7693
$ws.Range("L138").Value = 20419882.5
$ws.Range("M138").Value = 1222.2307
$ws.Range("N138").Value = -20430162.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2072.24
$ws.Range("I86").Value = 1984.1578
$ws.Range("J86").Value = 2351.1667
$ws.Range("K86").Value = 1984.1578
$ws.Range("L86").Value = 2351.1667
$ws.Range("M86").Value = -861.1578
$ws.Range("N86").Value = -4597.1667

$ws.Range("H89").Value = 2072.24
$ws.Range("I89").Value = 1984.1578
$ws.Range("J89").Value = 2351.1667
$ws.Range("K89").Value = 9920.789000000001
$ws.Range("L89").Value = 11755.8335
$ws.Range("M89").Value = -4304.789000000001
$ws.Range("N89").Value = -22987.8335

$ws.Range("H99").Value = 1416.625
$ws.Range("I99").Value = 1462.7273
$ws.Range("J99").Value = 1377.6154
$ws.Range("K99").Value = 1462.7273
$ws.Range("L99").Value = 1377.6154
$ws.Range("M99").Value = 35.27269999999999
$ws.Range("N99").Value = -4373.6154

$ws.Range("H107").Value = 1414.9166
$ws.Range("I107").Value = 997.1429000000001
$ws.Range("J107").Value = 1999.8
$ws.Range("K107").Value = 997.1429000000001
$ws.Range("L107").Value = 1999.8
$ws.Range("M107").Value = 922.8570999999999
$ws.Range("N107").Value = -5839.8

$ws.Range("H126").Value = 75000
$ws.Range("J126").Value = 75000
$ws.Range("L126").Value = 75000
$ws.Range("N126").Value = -84880

$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920

$ws.Range("H134").Value = 51228.332
$ws.Range("I134").Value = 3611.3125
$ws.Range("J134").Value = 203602.8
$ws.Range("K134").Value = 10833.9375
$ws.Range("L134").Value = 610808.3999999999
$ws.Range("M134").Value = -8298.9375
$ws.Range("N134").Value = -615878.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6603.3213
$ws.Range("I31").Value = 9891
$ws.Range("J31").Value = 3754
$ws.Range("K31").Value = 9891
$ws.Range("L31").Value = 3754
$ws.Range("M31").Value = -9596
$ws.Range("N31").Value = -4344

$ws.Range("H34").Value = 6603.3213
$ws.Range("I34").Value = 9891
$ws.Range("J34").Value = 3754
$ws.Range("K34").Value = 9891
$ws.Range("L34").Value = 3754
$ws.Range("M34").Value = -9689
$ws.Range("N34").Value = -4158

$ws.Range("H58").Value = 1492402.1
$ws.Range("I58").Value = 2393847.2
$ws.Range("J58").Value = 3057.8262
$ws.Range("K58").Value = 2393847.2
$ws.Range("L58").Value = 3057.8262
$ws.Range("M58").Value = -2393644.2
$ws.Range("N58").Value = -3463.8262

$ws.Range("H62").Value = 2418
$ws.Range("I62").Value = 2509.5833
$ws.Range("J62").Value = 2295.889
$ws.Range("K62").Value = 2509.5833
$ws.Range("L62").Value = 2295.889
$ws.Range("M62").Value = -1885.5833
$ws.Range("N62").Value = -3543.889

$ws.Range("H65").Value = 2418
$ws.Range("I65").Value = 2509.5833
$ws.Range("J65").Value = 2295.889
$ws.Range("K65").Value = 12547.9165
$ws.Range("L65").Value = 11479.445
$ws.Range("M65").Value = -9427.916499999999
$ws.Range("N65").Value = -17719.445

$ws.Range("H98").Value = 20750
$ws.Range("I98").Value = 15000
$ws.Range("J98").Value = 26500
$ws.Range("K98").Value = 15000
$ws.Range("L98").Value = 26500
$ws.Range("M98").Value = -12754
$ws.Range("N98").Value = -30992

$ws.Range("H100").Value = 63345
$ws.Range("J100").Value = 63345
$ws.Range("L100").Value = 63345
$ws.Range("N100").Value = -65509

$ws.Range("H103").Value = 35580
$ws.Range("I103").Value = 30725
$ws.Range("J103").Value = 55000
$ws.Range("K103").Value = 30725
$ws.Range("L103").Value = 55000
$ws.Range("M103").Value = -29553
$ws.Range("N103").Value = -57344

$ws.Range("H136").Value = 1492402.1
$ws.Range("I136").Value = 2393847.2
$ws.Range("J136").Value = 3057.8262
$ws.Range("K136").Value = 7181541.600000001
$ws.Range("L136").Value = 9173.4786
$ws.Range("M136").Value = -7178991.600000001
$ws.Range("N136").Value = -14273.4786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1865.7142
$ws.Range("I36").Value = 733.3333
$ws.Range("J36").Value = 2715
$ws.Range("K36").Value = 2199.9999
$ws.Range("L36").Value = 8145
$ws.Range("M36").Value = -2030.9999
$ws.Range("N36").Value = -8483

$ws.Range("H60").Value = 716.3077
$ws.Range("I60").Value = 415.83334
$ws.Range("J60").Value = 973.8570999999999
$ws.Range("K60").Value = 1247.50002
$ws.Range("L60").Value = 2921.5713
$ws.Range("M60").Value = -996.5000199999999
$ws.Range("N60").Value = -3423.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1572.1428
$ws.Range("I97").Value = 1656
$ws.Range("J97").Value = 1362.5
$ws.Range("K97").Value = 1656
$ws.Range("L97").Value = 1362.5
$ws.Range("M97").Value = -1160
$ws.Range("N97").Value = -2354.5

$ws.Range("H113").Value = 2484.6667
$ws.Range("I113").Value = 2562.1
$ws.Range("K113").Value = 2562.1
$ws.Range("M113").Value = -392.0999999999999

$ws.Range("H132").Value = 4784.6943
$ws.Range("I132").Value = 1956.4482
$ws.Range("J132").Value = 16501.715
$ws.Range("K132").Value = 5869.3446
$ws.Range("L132").Value = 49505.145
$ws.Range("M132").Value = -3339.3446
$ws.Range("N132").Value = -54565.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1125.1666
$ws.Range("I22").Value = 1050.2
$ws.Range("K22").Value = 1050.2
$ws.Range("M22").Value = -755.2

$ws.Range("H27").Value = 1125.1666
$ws.Range("I27").Value = 1050.2
$ws.Range("K27").Value = 1050.2
$ws.Range("M27").Value = -943.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3311.4546
$ws.Range("I62").Value = 3099.889
$ws.Range("K62").Value = 3099.889
$ws.Range("M62").Value = -2475.889

$ws.Range("H65").Value = 3311.4546
$ws.Range("I65").Value = 3099.889
$ws.Range("K65").Value = 15499.445
$ws.Range("M65").Value = -12379.445

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 1604.6383
$ws.Range("I132").Value = 652.36664
$ws.Range("K132").Value = 1957.09992
$ws.Range("M132").Value = 572.9000800000001
Write-Output "Applied all profit sheet updates"
